$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '68.111.71'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +7.69%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = "'" + '3.633.65'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +4.32%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').Value = "'" + '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.04%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = "'" + '421.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +1.69%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = "'" + '132.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +1.41%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').Value = "'" + '0.654'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +4.76%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').Value = "'" + '3.624.77'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +4.36%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').Value = "'" + '  -0.07%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Value = "'" + '0.772'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +5.68%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').Value = "'" + '0.198'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +32.29%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').Value = "'" + '0.0000426'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +93.56%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Value = "'" + '42.39'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -0.63%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = "'" + '9.83'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +1.10%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = "'" + '4.197.62'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +3.81%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('E16').Value = "'" + '  -0.21%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = "'" + '20.22'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -1.30%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = "'" + '3.624.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  +4.04%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Value = "'" + '1.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +2.76%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = "'" + '67.999.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +7.60%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').Value = "'" + '12.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -2.01%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = "'" + '458.75'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -1.76%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Value = "'" + '89.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -1.28%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').Value = "'" + '3.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -6.30%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').Value = "'" + '13.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +1.12%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Value = "'" + '3.34'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -1.58%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('B27').Value = "'" + 'Filecoin'
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').Value = "'" + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').Value = "'" + '10.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -3.13%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('B28').Value = "'" + 'EthereumClassic'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').Value = "'" + 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').Value = "'" + '35.97'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +6.61%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Value = "'" + '4.86'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +1.49%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').Value = "'" + '12.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +1.61%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Value = "'" + '2.77'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -1.55%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').Value = "'" + '0.120'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +6.64%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = "'" + '7.16'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  -5.25%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Value = "'" + '0.159'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -5.24%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').Value = "'" + '40.11'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -0.99%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').Value = "'" + '0.998'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +0.15%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').Value = "'" + '56.16'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -2.41%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').Value = "'" + '0.0486'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -0.26%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').Value = "'" + '0.0₃0760'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  +30.81%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').Value = "'" + '0.147'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +9.35%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').Value = "'" + '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -0.05%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').Value = "'" + 'Monero'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'" + 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'" + '147.85'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +0.03%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('B43').Value = "'" + 'Stacks'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'" + 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'" + '2.93'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -5.88%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = "'" + '2.66'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -6.23%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').Value = "'" + '3.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -2.57%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = "'" + '4.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -8.08%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('B47').Value = "'" + 'Cronos'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'" + 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'" + '0.168'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +18.44%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('B48').Value = "'" + 'TheGraph'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'" + 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'" + '0.303'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -4.31%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('B49').Value = "'" + 'ARBITRUM'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'" + 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'" + '1.99'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -1.90%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('B50').Value = "'" + 'ThetaToken'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'" + 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'" + '2.46'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  +4.65%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('B51').Value = "'" + 'ApeXProtocol'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'" + 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'" + '2.67'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +15.29%  '
$ws.Range('E51').Style = 'Normal'
